# Refresh cryptos worksheet: update price/volume columns (and reorder a few
# coin rows) to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.558.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.93%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.95%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.54"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.91%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4785"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.39%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.44%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07356"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.59%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9352"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.43%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.79"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07796"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.24%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.899.92"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.445"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.573"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.57%  "

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008832"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.99%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.012"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.30%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.606.41"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.87%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.69"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.105"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.72"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.66%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.935"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.25%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.19"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.15%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.51"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.27%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.026"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.88%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.49"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.81%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.957"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.36%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08881"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.12%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.333"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.209"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.48%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7572"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.58%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.598"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.92%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.695"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.37%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.126"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.09%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02036"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.17%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5679"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.63%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05361"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.93%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.979"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.33%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.050"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.27%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.528"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.84%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1528"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.67%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4894"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.27%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.66"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.08%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.37"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.31%  "

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.013"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.29%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.666"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.35%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "67.53"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.29%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06096"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9109"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.08%  "

